# COOLMATE-005 — define scenario: Add new address + set default (TC001),
# Update address (TC002), Delete address (TC003).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlContinuous = 1
$xlThin = 2
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# ---------------------------------------------------------------------
# 1. Row 2 — TC001: Add a new address and set default
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Add a new address to the address book with all valid information and set it as default"
$ws.Range("C2").Value = "A new address is created and set as default in the address book"
$ws.Range("D2").Value = "TS001-Address Book"
$ws.Range("E2").Value = "TC001-Add new address and set default"
$ws.Range("F2").Value = "Functional"
$ws.Range("G2").Value = "Auto"
$ws.Range("H2").Value = "N/A"
$ws.Rows.Item(2).RowHeight = 73.2

# ---------------------------------------------------------------------
# 2. Row 3 — TC002: Update address
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Update a address with all valid information in the address book"
$ws.Range("C3").Value = "The address is updated with valid information "
$ws.Range("D3").Value = "TS001-Address Book"
$ws.Range("E3").Value = "TC002-Update address"
$ws.Range("F3").Value = "Functional"
$ws.Range("G3").Value = "Auto"
$ws.Range("H3").Value = "N/A"
$ws.Rows.Item(3).RowHeight = 73.2

# ---------------------------------------------------------------------
# 3. Row 4 (new) — TC003: Delete address
#    Seed formatting from row 2 (same banding colour), then correct the
#    borders below.
# ---------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Delete an address in the address book"
$ws.Range("C4").Value = "The address is deleted"
$ws.Range("D4").Value = "TS001-Address Book"
$ws.Range("E4").Value = "TC-003-Delete address"
$ws.Range("F4").Value = "Functional"
$ws.Range("G4").Value = "Auto"
$ws.Range("H4").Value = "N/A"
$ws.Rows.Item(4).RowHeight = 59.25

# ---------------------------------------------------------------------
# 4. Borders — row 1 loses its bottom rule (row 2 below now owns the
#    line), row 2 becomes fully boxed (except col C, right-only), row 3
#    keeps only its vertical dividers, row 4 gets top+bottom+right rules
#    (left unboxed, matching the previous row-3 treatment) with column A
#    fully boxed.
# ---------------------------------------------------------------------

function Set-Edge($addr, $edge, $on) {
    $b = $ws.Range($addr).Borders.Item($edge)
    if ($on) {
        $b.LineStyle = $xlContinuous
        $b.Weight = $xlThin
    } else {
        $b.LineStyle = -4142
    }
}

# Row 1 — drop bottom border
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    Set-Edge "$($col)1" $xlEdgeBottom $false
}

# Row 2 — full box everywhere, column C stays right-only
foreach ($col in @("A","B","D","E","F","G","H")) {
    Set-Edge "$($col)2" $xlEdgeLeft $true
    Set-Edge "$($col)2" $xlEdgeTop $true
    Set-Edge "$($col)2" $xlEdgeBottom $true
    Set-Edge "$($col)2" $xlEdgeRight $true
}
Set-Edge "C2" $xlEdgeLeft $false
Set-Edge "C2" $xlEdgeTop $false
Set-Edge "C2" $xlEdgeBottom $false
Set-Edge "C2" $xlEdgeRight $true

# Row 3 — only vertical dividers (right edge); column A also keeps left edge
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    Set-Edge "$($col)3" $xlEdgeTop $false
    Set-Edge "$($col)3" $xlEdgeBottom $false
    Set-Edge "$($col)3" $xlEdgeRight $true
}
Set-Edge "A3" $xlEdgeLeft $true
foreach ($col in @("B","C","D","E","F","G","H")) {
    Set-Edge "$($col)3" $xlEdgeLeft $false
}

# Row 4 — top+bottom+right everywhere, column A also boxed on the left
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    Set-Edge "$($col)4" $xlEdgeTop $true
    Set-Edge "$($col)4" $xlEdgeBottom $true
    Set-Edge "$($col)4" $xlEdgeRight $true
}
Set-Edge "A4" $xlEdgeLeft $true
foreach ($col in @("B","C","D","E","F","G","H")) {
    Set-Edge "$($col)4" $xlEdgeLeft $false
}

# ---------------------------------------------------------------------
# 5. Sheet view — scroll back to column A, move the selection to C4
# ---------------------------------------------------------------------
[void]$ws.Range("C4").Select()
